$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora) is fully refreshed this run (row 2-51, "8" -> "9"); pre-format
# the whole column as Text so the numeric-looking value does not get auto-typed as a number.
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "335.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.99%"
$ws.Range("G2").Value = "9"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.94%"
$ws.Range("G3").Value = "9"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.802"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.86%"
$ws.Range("G4").Value = "9"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08341"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.14%"
$ws.Range("G5").Value = "9"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.813"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("G6").Value = "9"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.514"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.40%"
$ws.Range("G7").Value = "9"

# Row 8
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.986"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.51%"
$ws.Range("G8").Value = "9"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.898"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.72%"
$ws.Range("G9").Value = "9"

# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9410"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.24%"
$ws.Range("G10").Value = "9"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1235"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.37%"
$ws.Range("G11").Value = "9"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1957"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.12%"
$ws.Range("G12").Value = "9"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09881"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.69%"
$ws.Range("G13").Value = "9"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04546"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "19.43%"
$ws.Range("G14").Value = "9"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.86%"
$ws.Range("G15").Value = "9"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001314"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.05%"
$ws.Range("G16").Value = "9"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005964"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.13%"
$ws.Range("G17").Value = "9"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.496"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.39%"
$ws.Range("G18").Value = "9"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.74%"
$ws.Range("G19").Value = "9"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.789"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "7.82%"
$ws.Range("G20").Value = "9"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.68%"
$ws.Range("G21").Value = "9"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2610"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.22%"
$ws.Range("G22").Value = "9"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04399"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.05%"
$ws.Range("G23").Value = "9"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001258"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.02%"
$ws.Range("G24").Value = "9"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004391"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.46%"
$ws.Range("G25").Value = "9"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001261"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.92%"
$ws.Range("G26").Value = "9"

# Row 27
$ws.Range("G27").Value = "9"

# Row 28
$ws.Range("G28").Value = "9"

# Row 29
$ws.Range("G29").Value = "9"

# Row 30
$ws.Range("G30").Value = "9"

# Row 31
$ws.Range("G31").Value = "9"

# Row 32
$ws.Range("G32").Value = "9"

# Row 33
$ws.Range("G33").Value = "9"

# Row 34
$ws.Range("G34").Value = "9"

# Row 35
$ws.Range("G35").Value = "9"

# Row 36
$ws.Range("G36").Value = "9"

# Row 37
$ws.Range("G37").Value = "9"

# Row 38
$ws.Range("G38").Value = "9"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02794"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.82%"
$ws.Range("G39").Value = "9"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05711"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.58%"
$ws.Range("G40").Value = "9"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007921"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.34%"
$ws.Range("G41").Value = "9"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.97%"
$ws.Range("G42").Value = "9"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008969"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.21%"
$ws.Range("G43").Value = "9"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.30%"
$ws.Range("G44").Value = "9"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01059"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-7.91%"
$ws.Range("G45").Value = "9"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007287"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "10.46%"
$ws.Range("G46").Value = "9"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("G47").Value = "9"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003243"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.20%"
$ws.Range("G48").Value = "9"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.47%"
$ws.Range("G49").Value = "9"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("G50").Value = "9"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
$ws.Range("G51").Value = "9"
